$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 58:65, column F ("수집일자"/collection-date column) had been left holding
# raw date serials (45636) formatted with a yyyy-mm-dd number format, while every
# other row (2:57) stores the same information as the literal text "2024-12-09".
# Fix F58:F65 so they match the rest of the column: plain text "2024-12-09" with
# the same (unbordered, default) style used by F2:F57.

# Mark the cells as text first so typing a date-shaped string doesn't get
# reinterpreted by Excel's smart entry as a date serial again.
$ws.Range("F58:F65").NumberFormat = "@"
$ws.Range("F58:F65").Value = "2024-12-09"

# Now copy the formatting (font/number format/borders) from F2 - a cell that
# already has the desired look - onto F58:F65 so the style matches the rest
# of the column exactly.
$ws.Range("F2").Copy()
$ws.Range("F58:F65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen column F so the "2024-12-09" text is fully visible.
$ws.Columns("F").ColumnWidth = 19.3

# Leave the selection where the edit was made.
$ws.Range("F65").Select()
